$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0, "roll_crawl", 0.55, 1.2, 1.3),
    @(2, "roll_crawl", 0.5499999999999999, 0, 1.3),
    @(3, "roll_crawl", 0.5333333333333333, 0, 1.3),
    @(4, "roll_crawl", 0.55, 0, 1.3),
    @(5, "roll_roll", 0.54, 0, 1.3)
)

$row = 6
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $ws.Cells.Item($row, 4).Value = $r[3]
    $ws.Cells.Item($row, 5).Value = $r[4]
    $row++
}
